# here_1.0.1 traceability-matrix refresh: new title/package/date header row,
# exported_function header row moved down, per-function descriptions added in a
# wide column D, and the now-unused Signature column (F) is dropped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old "Signature" column (F) entirely - no longer part of the sheet.
$ws.Columns.Item(6).Delete() | Out-Null

# Re-home the header formatting (italic font / purple fill / centered / bottom
# border) from row 1 (old header) down to row 2 (new header), since row 1 is
# now a title/metadata row with default formatting.
$ws.Range("A1:E1").Copy() | Out-Null
$ws.Range("A2:E2").PasteSpecial(-4122) | Out-Null
$ws.Range("A1:E1").ClearFormats() | Out-Null

# Write the new cell contents.
$ws.Range("A1").Value = 'Traceability Matrix'
$ws.Range("B1").Value = 'Package'
$ws.Range("C1").Value = 'here'
$ws.Range("D1").Value = 'Date Time'
$ws.Range("E1").Value = '2024-06-07 12:26:36'
$ws.Range("A2").Value = 'exported_function'
$ws.Range("B2").Value = 'code_script'
$ws.Range("C2").Value = 'documentation'
$ws.Range("D2").Value = 'description'
$ws.Range("E2").Value = 'coverage_percent'
$ws.Range("A3").Value = 'dr_here'
$ws.Range("B3").Value = 'R/dr_here.R'
$ws.Range("C3").Value = 'dr_here.Rd'
$ws.Range("D3").Value = 'dr_here() shows a message that by default also includes thereason why here() is set to a particular directory.Use this function if here() gives unexpected results.'
$ws.Range("E3").Value = 100
$ws.Range("A4").Value = 'here'
$ws.Range("B4").Value = 'R/here.R'
$ws.Range("C4").Value = 'here.Rd'
$ws.Range("D4").Value = 'here() uses a reasonable heuristics to find your project''s files, based onthe current working directory at the time when the package is loaded.Use it as a drop-in replacement for file.path(), it will always locate thefiles relative to your project root.'
$ws.Range("E4").Value = 100
$ws.Range("A5").Value = 'i_am'
$ws.Range("B5").Value = 'R/i_am.R'
$ws.Range("C5").Value = 'i_am.Rd'
$ws.Range("D5").Value = 'Add a call to here::i_am("<project-relative path>.<ext>")at the top of your R script or in the first chunk of your rmarkdown document.This ensures that the project root is set up correctly:subsequent calls to here() will refer to the implied project root.If the current working directory is outside of the projectwhere the script or report is intended to run, it will failwith a descriptive message.'
$ws.Range("E5").Value = 95.83
$ws.Range("A6").Value = 'set_here'
$ws.Range("B6").Value = 'R/set_here.R'
$ws.Range("C6").Value = 'set_here.Rd'
$ws.Range("D6").Value = 'html<a href=''https://www.tidyverse.org/lifecycle/#superseded''><img src=''figures/lifecycle-superseded.svg'' alt=''Superseded lifecycle''></a>Superseded'
$ws.Range("E6").Value = 100

# Resize columns: A-C modest widths, D very wide (holds long descriptions), E narrower.
$ws.Columns.Item(1).ColumnWidth = 18.76
$ws.Columns.Item(2).ColumnWidth = 13.76
$ws.Columns.Item(3).ColumnWidth = 14.76
$ws.Columns.Item(4).ColumnWidth = 400.76
$ws.Columns.Item(5).ColumnWidth = 17.76
